$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings (e.g. "1.002")
# are preserved verbatim as text instead of being parsed into numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = '31.054.46'
$ws.Range("E2").Value = '  +1.52%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.957.43'
$ws.Range("E3").Value = '  +1.84%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.28%  '

# Row 5: BNB
$ws.Range("D5").Value = '246.54'
$ws.Range("E5").Value = '  +0.50%  '

# Row 6: USDC
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.19%  '

# Row 7: XRP
$ws.Range("D7").Value = '0.4906'
$ws.Range("E7").Value = '  +1.67%  '

# Row 8: OKB
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '44.52'
$ws.Range("E8").Value = '  -0.09%  '

# Row 9: Cardano
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.2965'
$ws.Range("E9").Value = '  +2.25%  '

# Row 10: Dogecoin
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.06815'
$ws.Range("E10").Value = '  +0.07%  '

# Row 11: Solana
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '18.99'
$ws.Range("E11").Value = '  -2.43%  '

# Row 12: Litecoin
$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").Value = '106.25'
$ws.Range("E12").Value = '  -5.32%  '

# Row 13: TRON
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.07744'
$ws.Range("E13").Value = '  +2.23%  '

# Row 14: WrappedEther
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.931.43'
$ws.Range("E14").Value = '  +0.55%  '

# Row 15: Polkadot
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '5.406'
$ws.Range("E15").Value = '  -1.66%  '

# Row 16: Polygon
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = '0.7090'
$ws.Range("E16").Value = '  +5.21%  '

# Row 17: BitcoinCash
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").Value = '283.23'
$ws.Range("E17").Value = '  -4.07%  '

# Row 18: WrappedBTC
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '31.092.89'
$ws.Range("E18").Value = '  +1.70%  '

# Row 19: ShibaInu
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000007767'
$ws.Range("E19").Value = '  +1.15%  '

# Row 20: Avalanche
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '13.23'
$ws.Range("E20").Value = '  +1.44%  '

# Row 21: Dai
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.17%  '

# Row 22: WrappedliquidstakedEther2.0
$ws.Range("D22").Value = '2.187.74'
$ws.Range("E22").Value = '  +1.12%  '

# Row 23: Uniswap
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '5.549'
$ws.Range("E23").Value = '  +0.73%  '

# Row 24: BinanceUSD
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  +0.22%  '

# Row 25: Chainlink
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").Value = '6.590'
$ws.Range("E25").Value = '  +2.14%  '

# Row 26: Cosmos
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '9.964'
$ws.Range("E26").Value = '  +4.96%  '

# Row 27: Monero
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '168.55'
$ws.Range("E27").Value = '  +0.81%  '

# Row 28: EthereumClassic
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '19.97'
$ws.Range("E28").Value = '  -1.94%  '

# Row 29: LidoDAOToken
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '2.192'
$ws.Range("E29").Value = '  +4.76%  '

# Row 30: Stellar
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.1060'
$ws.Range("E30").Value = '  -0.55%  '

# Row 31: Toncoin
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '1.442'
$ws.Range("E31").Value = '  -0.02%  '

# Row 32: Filecoin
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '4.781'
$ws.Range("E32").Value = '  +17.68%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '4.524'
$ws.Range("E33").Value = '  +9.35%  '

# Row 34: Hedera
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.05003'
$ws.Range("E34").Value = '  +0.41%  '

# Row 35: ImmutableX
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.7658'
$ws.Range("E35").Value = '  +4.24%  '

# Row 36: ARBITRUM
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.165'
$ws.Range("E36").Value = '  +2.35%  '

# Row 37: VeChain
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02052'
$ws.Range("E37").Value = '  +1.08%  '

# Row 38: HuobiToken
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").Value = '2.733'
$ws.Range("E38").Value = '  +0.69%  '

# Row 39: MXToken
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.713'
$ws.Range("E39").Value = '  +1.05%  '

# Row 40: RenderToken
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '2.130'
$ws.Range("E40").Value = '  +5.16%  '

# Row 41: FraxShare
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '6.416'
$ws.Range("E41").Value = '  +9.32%  '

# Row 42: TrustWalletToken
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '0.8851'
$ws.Range("E42").Value = '  +1.72%  '

# Row 43: Quant
$ws.Range("D43").Value = '109.40'
$ws.Range("E43").Value = '  -0.10%  '

# Row 44: TheSandbox
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.4457'
$ws.Range("E44").Value = '  +0.46%  '

# Row 45: Aave
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '73.11'
$ws.Range("E45").Value = '  +5.15%  '

# Row 46: PaxDollar
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.01%  '

# Row 47: Maker
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.021.29'
$ws.Range("E47").Value = '  +21.18%  '

# Row 48: Aptos
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.468'
$ws.Range("E48").Value = '  +2.85%  '

# Row 49: Algorand
$ws.Range("D49").Value = '0.1267'
$ws.Range("E49").Value = '  +3.12%  '

# Row 50: EnergySwap
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.371'
$ws.Range("E50").Value = '  +1.55%  '

# Row 51: Elrond
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '35.96'
$ws.Range("E51").Value = '  +3.11%  '

# Restore General format on the price column now that the literal text is committed,
# so no stray quote-prefix/text style lingers on these cells.
$priceRange.ClearFormats()
